$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New URL data replacing the old columbiachile.cl data
$ws.Range("A2").Value = "https://www.zapatos.cl/"
$ws.Range("A3").Value = "https://www.zapatos.cl/?__disablePixels"
$ws.Range("A4").Value = "https://www.zapatos.cl/mujer?map=genero"
$ws.Range("A5").Value = "https://www.zapatos.cl/mujer?map=genero?__disablePixels"
$ws.Range("A6").Value = "https://www.zapatos.cl/polera-m-c-mujer-v-neck-bsoul-bs210021428-7un/p"
$ws.Range("A7").Value = "https://www.zapatos.cl/polera-m-c-mujer-v-neck-bsoul-bs210021428-7un/p?__disablePixels"
$ws.Range("A8").Value = "https://www.zapatos.cl/checkout/?orderFormId=f90f4001500640a9a603a60e61c61d60#/cart"

# Apply the same "Hipervinculo" style to all data rows (A2 included) so they
# all share one cell format (General number format, not Text)
$ws.Range("A2:A8").Style = "Hipervínculo"

# Update selection to match the new extent of data
$ws.Range("A2:A8").Select()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
